# Scheduled market-data refresh: update cached price/profit figures
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# for the affected leves across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 315.85715
$ws.Range("I18").Value = 321.83334
$ws.Range("J18").Value = 280
$ws.Range("K18").Value = 321.83334
$ws.Range("L18").Value = 280
$ws.Range("M18").Value = -37.83334000000002
$ws.Range("N18").Value = -848

$ws.Range("H51").Value = 1809.1818
$ws.Range("I51").Value = 5750.5
$ws.Range("J51").Value = 933.3333
$ws.Range("K51").Value = 5750.5
$ws.Range("L51").Value = 933.3333
$ws.Range("M51").Value = -5266.5
$ws.Range("N51").Value = -1901.3333

$ws.Range("H64").Value = 3502.1724
$ws.Range("I64").Value = 3302.9092
$ws.Range("J64").Value = 4128.4287
$ws.Range("K64").Value = 3302.9092
$ws.Range("L64").Value = 4128.4287
$ws.Range("M64").Value = -3054.9092
$ws.Range("N64").Value = -4624.4287

$ws.Range("H67").Value = 3502.1724
$ws.Range("I67").Value = 3302.9092
$ws.Range("J67").Value = 4128.4287
$ws.Range("K67").Value = 3302.9092
$ws.Range("L67").Value = 4128.4287
$ws.Range("M67").Value = -2444.9092
$ws.Range("N67").Value = -5844.4287

$ws.Range("H76").Value = 3572.325
$ws.Range("I76").Value = 3474.7812
$ws.Range("K76").Value = 3474.7812
$ws.Range("M76").Value = -3159.7812

$ws.Range("H79").Value = 3572.325
$ws.Range("I79").Value = 3474.7812
$ws.Range("K79").Value = 3474.7812
$ws.Range("M79").Value = -2382.7812

$ws.Range("H111").Value = 5225.875
$ws.Range("I111").Value = 6961.4
$ws.Range("K111").Value = 20884.2
$ws.Range("M111").Value = -17817.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 2500
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -3872

$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 12500
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -19364

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 3750
$ws.Range("I42").Value = 3750
$ws.Range("K42").Value = 3750
$ws.Range("M42").Value = -3157

$ws.Range("H62").Value = 3125
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 3500
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 3500
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -4748

$ws.Range("H65").Value = 3125
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 3500
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 17500
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -23740

$ws.Range("H107").Value = 392.75
$ws.Range("I107").Value = 387
$ws.Range("J107").Value = 421.5
$ws.Range("K107").Value = 387
$ws.Range("L107").Value = 421.5
$ws.Range("M107").Value = 1533
$ws.Range("N107").Value = -4261.5

$ws.Range("H132").Value = 2928.238
$ws.Range("I132").Value = 2614
$ws.Range("J132").Value = 3347.2222
$ws.Range("K132").Value = 7842
$ws.Range("L132").Value = 10041.6666
$ws.Range("M132").Value = -5312
$ws.Range("N132").Value = -15101.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 12503038
$ws.Range("I5").Value = 354.6
$ws.Range("J5").Value = 50011090
$ws.Range("K5").Value = 1063.8
$ws.Range("L5").Value = 150033270
$ws.Range("M5").Value = -951.8000000000002
$ws.Range("N5").Value = -150033494

$ws.Range("H97").Value = 10006.083
$ws.Range("J97").Value = 16467.143
$ws.Range("L97").Value = 49401.429
$ws.Range("N97").Value = -50393.429

$ws.Range("H113").Value = 712.84
$ws.Range("I113").Value = 740.7975
$ws.Range("J113").Value = 607.6667
$ws.Range("K113").Value = 2222.3925
$ws.Range("L113").Value = 1823.0001
$ws.Range("M113").Value = -52.39249999999993
$ws.Range("N113").Value = -6163.0001

$ws.Range("H121").Value = 220
$ws.Range("J121").Value = 200
$ws.Range("L121").Value = 600
$ws.Range("N121").Value = -3220

$ws.Range("H122").Value = 537.8333
$ws.Range("I122").Value = 437.53333
$ws.Range("J122").Value = 638.13336
$ws.Range("K122").Value = 3937.79997
$ws.Range("L122").Value = 5743.20024
$ws.Range("M122").Value = -1487.79997
$ws.Range("N122").Value = -10643.20024

$ws.Range("H129").Value = 2423.8
$ws.Range("I129").Value = 2807.5
$ws.Range("J129").Value = 2284.2727
$ws.Range("K129").Value = 8422.5
$ws.Range("L129").Value = 6852.8181
$ws.Range("M129").Value = -3422.5
$ws.Range("N129").Value = -16852.8181

$ws.Range("H131").Value = 1105.5264
$ws.Range("I131").Value = 1413.2142
$ws.Range("J131").Value = 926.0417
$ws.Range("K131").Value = 4239.642599999999
$ws.Range("L131").Value = 2778.1251
$ws.Range("M131").Value = 800.3574000000008
$ws.Range("N131").Value = -12858.1251

$ws.Range("H135").Value = 12503038
$ws.Range("I135").Value = 354.6
$ws.Range("J135").Value = 50011090
$ws.Range("K135").Value = 3191.4
$ws.Range("L135").Value = 450099810
$ws.Range("M135").Value = -656.4000000000001
$ws.Range("N135").Value = -450104880

$ws.Range("H139").Value = 2073889.8
$ws.Range("I139").Value = 5420135
$ws.Range("J139").Value = 2404.476
$ws.Range("K139").Value = 16260405
$ws.Range("L139").Value = 7213.428
$ws.Range("M139").Value = -16255265
$ws.Range("N139").Value = -17493.428

$ws.Range("H140").Value = 2224.6667
$ws.Range("I140").Value = 1792.5
$ws.Range("J140").Value = 3089
$ws.Range("K140").Value = 5377.5
$ws.Range("L140").Value = 9267
$ws.Range("M140").Value = -197.5
$ws.Range("N140").Value = -19627

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4914.385
$ws.Range("I70").Value = 4441.5625
$ws.Range("J70").Value = 5670.9
$ws.Range("K70").Value = 4441.5625
$ws.Range("L70").Value = 5670.9
$ws.Range("M70").Value = -4171.5625
$ws.Range("N70").Value = -6210.9

$ws.Range("H73").Value = 4914.385
$ws.Range("I73").Value = 4441.5625
$ws.Range("J73").Value = 5670.9
$ws.Range("K73").Value = 4441.5625
$ws.Range("L73").Value = 5670.9
$ws.Range("M73").Value = -3505.5625
$ws.Range("N73").Value = -7542.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5070.104
$ws.Range("I136").Value = 2614.348
$ws.Range("J136").Value = 7329.4
$ws.Range("K136").Value = 7843.044
$ws.Range("L136").Value = 21988.2
$ws.Range("M136").Value = -5293.044
$ws.Range("N136").Value = -27088.2
